# Applies the commit:
# "fixed 2 problems: missed a few spots to multiple eggs by female prop
#  and re-ran with corrected inital condition"
#
# Concretely, on sheet "7_spp_GN_enmalle" a new row is inserted right
# below the "Corvina reina" data row (old row 6). The new row converts
# that row's raw catch numbers into corrected figures using the female
# proportion factor kept in ezfureza_por_arte_y_spp_2012!$K$9, pushing
# every row from the old row 7 onward down by one.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("7_spp_GN_enmalle")

# Insert a new blank row 7 (formats/styles are inherited from row 6,
# same as Excel's normal "Insert Row" behavior).
$ws1.Rows.Item(7).Insert()

# Fill in the corrected values as formulas, referencing the row above
# (the raw "Corvina reina" catch numbers) and the correction factor on
# the ezfureza_por_arte_y_spp_2012 sheet.
$ws1.Range("C7").Formula = "=C6*1000/ezfureza_por_arte_y_spp_2012!`$K`$9"
$ws1.Range("D7").Formula = "=D6*1000/ezfureza_por_arte_y_spp_2012!`$K`$9"
$ws1.Range("E7").Formula = "=E6*1000/ezfureza_por_arte_y_spp_2012!`$K`$9"
$ws1.Range("F7").Formula = "=F6*1000/ezfureza_por_arte_y_spp_2012!`$K`$9"
$ws1.Range("G7").Formula = "=G6*1000/ezfureza_por_arte_y_spp_2012!`$K`$9"

# Move the active selection/tab to sheet "7_spp_GN_enmalle" at F22 (this
# also clears the tabSelected flag that used to sit on
# ezfureza_por_arte_y_spp_2012).
$ws1.Activate()
$ws1.Range("F22").Select() | Out-Null
